$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45178 (2023-09-09)
# to 45179 (2023-09-10) for every data row (rows 2 through 23).
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
